# Scheduled-runner data refresh: update the per-Leve market-board price/profit
# columns (H:N) on each job sheet to reflect the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1174.375
$ws.Range("I62").Value = 1159
$ws.Range("J62").Value = 1200
$ws.Range("K62").Value = 1159
$ws.Range("L62").Value = 1200
$ws.Range("M62").Value = -535
$ws.Range("N62").Value = -2448
$ws.Range("H65").Value = 1174.375
$ws.Range("I65").Value = 1159
$ws.Range("J65").Value = 1200
$ws.Range("K65").Value = 5795
$ws.Range("L65").Value = 6000
$ws.Range("M65").Value = -2675
$ws.Range("N65").Value = -12240
$ws.Range("H137").Value = 9524716
$ws.Range("I137").Value = 918.25
$ws.Range("J137").Value = 40000868
$ws.Range("K137").Value = 2754.75
$ws.Range("L137").Value = 120002604
$ws.Range("M137").Value = -204.75
$ws.Range("N137").Value = -120007704
$ws.Range("H138").Value = 1812.4286
$ws.Range("I138").Value = 1717.2307
$ws.Range("J138").Value = 3050
$ws.Range("K138").Value = 5151.6921
$ws.Range("L138").Value = 9150
$ws.Range("M138").Value = -11.69210000000021
$ws.Range("N138").Value = -19430
$ws.Range("H141").Value = 877.41815
$ws.Range("I141").Value = 810.6226
$ws.Range("J141").Value = 2647.5
$ws.Range("K141").Value = 2431.8678
$ws.Range("L141").Value = 7942.5
$ws.Range("M141").Value = 2748.1322
$ws.Range("N141").Value = -18302.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 125096.625
$ws.Range("I5").Value = 166762.17
$ws.Range("K5").Value = 166762.17
$ws.Range("M5").Value = -166650.17
$ws.Range("H45").Value = 1627.8077
$ws.Range("I45").Value = 1729.2858
$ws.Range("J45").Value = 1509.4166
$ws.Range("K45").Value = 1729.2858
$ws.Range("L45").Value = 1509.4166
$ws.Range("M45").Value = -1352.2858
$ws.Range("N45").Value = -2263.4166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 125096.625
$ws.Range("I4").Value = 166762.17
$ws.Range("K4").Value = 166762.17
$ws.Range("M4").Value = -166647.17
$ws.Range("H80").Value = 13783.467
$ws.Range("I80").Value = 40316
$ws.Range("J80").Value = 517.2
$ws.Range("K80").Value = 40316
$ws.Range("L80").Value = 517.2
$ws.Range("M80").Value = -39318
$ws.Range("N80").Value = -2513.2
$ws.Range("H83").Value = 13783.467
$ws.Range("I83").Value = 40316
$ws.Range("J83").Value = 517.2
$ws.Range("K83").Value = 201580
$ws.Range("L83").Value = 2586
$ws.Range("M83").Value = -196588
$ws.Range("N83").Value = -12570
$ws.Range("H105").Value = 4331.1055
$ws.Range("I105").Value = 3002
$ws.Range("J105").Value = 4805.7856
$ws.Range("K105").Value = 3002
$ws.Range("L105").Value = 4805.7856
$ws.Range("M105").Value = -1255
$ws.Range("N105").Value = -8299.785599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 122.14286
$ws.Range("I7").Value = 150.33333
$ws.Range("J7").Value = 101
$ws.Range("K7").Value = 150.33333
$ws.Range("L7").Value = 101
$ws.Range("M7").Value = -37.33332999999999
$ws.Range("N7").Value = -327
$ws.Range("H132").Value = 6494750.5
$ws.Range("I132").Value = 9434937
$ws.Range("J132").Value = 1837.8334
$ws.Range("K132").Value = 28304811
$ws.Range("L132").Value = 5513.5002
$ws.Range("M132").Value = -28302281
$ws.Range("N132").Value = -10573.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H122").Value = 783.53845
$ws.Range("I122").Value = 726.2727
$ws.Range("J122").Value = 1098.5
$ws.Range("K122").Value = 6536.454299999999
$ws.Range("L122").Value = 9886.5
$ws.Range("M122").Value = -4086.454299999999
$ws.Range("N122").Value = -14786.5
$ws.Range("H136").Value = 2386.0588
$ws.Range("I136").Value = 1383
$ws.Range("J136").Value = 3819
$ws.Range("K136").Value = 4149
$ws.Range("L136").Value = 11457
$ws.Range("M136").Value = 951
$ws.Range("N136").Value = -21657
$ws.Range("H138").Value = 4336.1763
$ws.Range("I138").Value = 902.4
$ws.Range("K138").Value = 2707.2
$ws.Range("M138").Value = 2432.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 58.666668
$ws.Range("I2").Value = 62
$ws.Range("J2").Value = 53.42857
$ws.Range("K2").Value = 62
$ws.Range("L2").Value = 53.42857
$ws.Range("M2").Value = 51
$ws.Range("N2").Value = -279.42857
$ws.Range("H15").Value = 15036.286
$ws.Range("J15").Value = 15036.286
$ws.Range("L15").Value = 15036.286
$ws.Range("N15").Value = -15612.286
$ws.Range("H81").Value = 15036.286
$ws.Range("J81").Value = 15036.286
$ws.Range("L81").Value = 15036.286
$ws.Range("N81").Value = -17032.286
$ws.Range("H82").Value = 38000
$ws.Range("J82").Value = 38000
$ws.Range("L82").Value = 38000
$ws.Range("N82").Value = -38766
$ws.Range("H84").Value = 15036.286
$ws.Range("J84").Value = 15036.286
$ws.Range("L84").Value = 45108.858
$ws.Range("N84").Value = -55092.858
$ws.Range("H85").Value = 38000
$ws.Range("J85").Value = 38000
$ws.Range("L85").Value = 38000
$ws.Range("N85").Value = -40652
$ws.Range("H102").Value = 6953.6
$ws.Range("I102").Value = 7170.6665
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 7170.6665
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -5548.6665
$ws.Range("N102").Value = -8244
$ws.Range("H126").Value = 5464.647
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5685.8857
$ws.Range("J7").Value = 5085.5
$ws.Range("L7").Value = 5085.5
$ws.Range("N7").Value = -5309.5
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2830
$ws.Range("H81").Value = 37911.75
$ws.Range("J81").Value = 37911.75
$ws.Range("L81").Value = 37911.75
$ws.Range("N81").Value = -39907.75
$ws.Range("H84").Value = 37911.75
$ws.Range("J84").Value = 37911.75
$ws.Range("L84").Value = 113735.25
$ws.Range("N84").Value = -123719.25
$ws.Range("H126").Value = 5685.8857
$ws.Range("J126").Value = 5085.5
$ws.Range("L126").Value = 15256.5
$ws.Range("N126").Value = -20196.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 34886
$ws.Range("J75").Value = 34886
$ws.Range("L75").Value = 34886
$ws.Range("N75").Value = -36758
$ws.Range("H78").Value = 34886
$ws.Range("J78").Value = 34886
$ws.Range("L78").Value = 104658
$ws.Range("N78").Value = -114018
$ws.Range("H80").Value = 39180.6
$ws.Range("I80").Value = 40000
$ws.Range("J80").Value = 38975.75
$ws.Range("K80").Value = 40000
$ws.Range("L80").Value = 38975.75
$ws.Range("M80").Value = -39002
$ws.Range("N80").Value = -40971.75
$ws.Range("H83").Value = 39180.6
$ws.Range("I83").Value = 40000
$ws.Range("J83").Value = 38975.75
$ws.Range("K83").Value = 120000
$ws.Range("L83").Value = 116927.25
$ws.Range("M83").Value = -115008
$ws.Range("N83").Value = -126911.25
$ws.Range("H113").Value = 1132.1578
$ws.Range("I113").Value = 395.75
$ws.Range("J113").Value = 2394.5715
$ws.Range("K113").Value = 1187.25
$ws.Range("L113").Value = 7183.7145
$ws.Range("M113").Value = 982.75
$ws.Range("N113").Value = -11523.7145
